# "Generate Report for Archive"
#
# The localization-status report is regenerated: every cell that showed the
# placeholder status text "Ready for handoff" now reads "In Translation"
# (Overview!E2:F3 mirrors the zh-cn/de-de sheets' Status column, and each
# language sheet's own Status column, zh-cn!C2:C3 / de-de!C2:C3).
#
# The Status columns are also narrower afterwards (their width had been
# sized to fit the old, longer "Ready for handoff" text; after the text
# shrinks to "In Translation" the report-generator re-sizes them down).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
